$wb = $excel.ActiveWorkbook

$ws4 = $wb.Worksheets.Item("Mix_Mass_Frac")
$ws4.Range("B2").Value = 0.4174279927961577
$ws4.Range("C2").Value = 0.07594363404940722
$ws4.Range("D2").Value = 0.004219131060181059
$ws4.Range("E2").Value = 0.11461536284168035
$ws4.Range("F2").Value = 0.006909000270921163
$ws4.Range("G2").Value = 0.021785783678512083
$ws4.Range("I2").Value = 0.287732562868526
$ws4.Range("J2").Value = 0.4242786391451718
$ws4.Range("L2").Value = 0.0007133914144382479
$ws4.Range("N2").Value = 0.06380249467116196
$ws4.Range("B3").Value = 0.3739312249795685
$ws4.Range("C3").Value = 0.08300468581932795
$ws4.Range("D3").Value = 0.004613479262540833
$ws4.Range("E3").Value = 0.12930083699642483
$ws4.Range("F3").Value = 0.00734231521538458
$ws4.Range("G3").Value = 0.019997835516329554
$ws4.Range("I3").Value = 0.226179314894383
$ws4.Range("J3").Value = 0.4538290429783439
$ws4.Range("L3").Value = 0.0008135548307047094
$ws4.Range("N3").Value = 0.07491893448656062
$ws4.Range("B4").Value = 0.24400982245409192
$ws4.Range("C4").Value = 0.09522988090519836
$ws4.Range("D4").Value = 0.003549235559306227
$ws4.Range("E4").Value = 0.12321852723611326
$ws4.Range("F4").Value = 0.00908060514404102
$ws4.Range("G4").Value = 0.041455826886504184
$ws4.Range("I4").Value = 0.099510439123668
$ws4.Range("J4").Value = 0.5626909810288863
$ws4.Range("L4").Value = 0.0006043546885204292
$ws4.Range("N4").Value = 0.06466014942776233
$ws4.Range("B5").Value = 0.173716620121865
$ws4.Range("C5").Value = 0.09631250773927959
$ws4.Range("D5").Value = 0.0025060546242810615
$ws4.Range("E5").Value = 0.11793103702537136
$ws4.Range("F5").Value = 0.009968901808635522
$ws4.Range("G5").Value = 0.06162055808232024
$ws4.Range("I5").Value = 0.04886898051942753
$ws4.Range("J5").Value = 0.6152445286248649
$ws4.Range("L5").Value = 0.0004008017742054667
$ws4.Range("N5").Value = 0.04714662980161428
$ws4.Range("B6").Value = 0.12592825492062953
$ws4.Range("C6").Value = 0.0975568724384947
$ws4.Range("D6").Value = 0.0013095539678176804
$ws4.Range("E6").Value = 0.09974227591637293
$ws4.Range("F6").Value = 0.010404302189858937
$ws4.Range("G6").Value = 0.09662459453837795
$ws4.Range("I6").Value = 0.025268929829495763
$ws4.Range("J6").Value = 0.6433463342058855
$ws4.Range("L6").Value = 0.0002164445962176009
$ws4.Range("N6").Value = 0.0255306923174788
$ws4.Range("B7").Value = 0.08371857778173938
$ws4.Range("C7").Value = 0.08229902872160454
$ws4.Range("D7").Value = 0.0008054055992490316
$ws4.Range("E7").Value = 0.08077073533005223
$ws4.Range("F7").Value = 0.011181294686439296
$ws4.Range("G7").Value = 0.1171752365211685
$ws4.Range("I7").Value = 0.002995404528409153
$ws4.Range("J7").Value = 0.6866803637307566
$ws4.Range("L7").Value = 0.00013179098418896262
$ws4.Range("N7").Value = 0.01796073989813185
$ws4.Range("B8").Value = 0.04650516819806561
$ws4.Range("C8").Value = 0.05841305279482978
$ws4.Range("D8").Value = 0.00015870375684867312
$ws4.Range("E8").Value = 0.04803411772434896
$ws4.Range("F8").Value = 0.011796185776060016
$ws4.Range("G8").Value = 0.15422841289381825
$ws4.Range("I8").Value = 0.000318031722244082
$ws4.Range("J8").Value = 0.7238678631559408
$ws4.Range("L8").Value = 0.000023188282481681027
$ws4.Range("N8").Value = 0.003160443893427675
$ws4.Range("B9").Value = 0.030063208332282338
$ws4.Range("C9").Value = 0.039495533522106316
$ws4.Range("E9").Value = 0.03269526254054541
$ws4.Range("F9").Value = 0.011855314193272237
$ws4.Range("G9").Value = 0.17995647932240258
$ws4.Range("I9").Value = 0.0005778596112960189
$ws4.Range("J9").Value = 0.7347829650987584
$ws4.Range("L9").Value = 0.00000442616180071564
$ws4.Range("N9").Value = 0.0006321595498183332
$ws4.Range("B10").Value = 0.010527033763173078
$ws4.Range("C10").Value = 0.014459886264583246
$ws4.Range("E10").Value = 0.011865009655018815
$ws4.Range("F10").Value = 0.012145176273342695
$ws4.Range("G10").Value = 0.21306943527723599
$ws4.Range("J10").Value = 0.7484604925298192

$ws5 = $wb.Worksheets.Item("Uncertainties")
$ws5.Range("B2").Value = 0.07299117691602462
$ws5.Range("C2").Value = 0.009205011706618246
$ws5.Range("D2").Value = 0.00033609724666102563
$ws5.Range("E2").Value = 0.022979121317208685
$ws5.Range("F2").Value = 0.0006137043244842242
$ws5.Range("G2").Value = 0.0020569928391430204
$ws5.Range("I2").Value = 0.0725361739161283
$ws5.Range("J2").Value = 0.026771951190910557
$ws5.Range("L2").Value = 0.00006299542752974483
$ws5.Range("N2").Value = 0.004033820423441969
$ws5.Range("B3").Value = 0.11859774332441761
$ws5.Range("C3").Value = 0.011204041524548663
$ws5.Range("D3").Value = 0.0006768792550240583
$ws5.Range("E3").Value = 0.033499305283460984
$ws5.Range("F3").Value = 0.0012097500210485378
$ws5.Range("G3").Value = 0.011758540626339927
$ws5.Range("I3").Value = 0.1177652160048267
$ws5.Range("J3").Value = 0.06385795069971287
$ws5.Range("L3").Value = 0.00013663296247674346
$ws5.Range("N3").Value = 0.009974228791123675
$ws5.Range("B4").Value = 0.03989012605672048
$ws5.Range("C4").Value = 0.010224437265368716
$ws5.Range("D4").Value = 0.0005095689670170224
$ws5.Range("E4").Value = 0.03687671521963614
$ws5.Range("F4").Value = 0.0010677610586162034
$ws5.Range("G4").Value = 0.010389077804728837
$ws5.Range("I4").Value = 0.03817027533118383
$ws5.Range("J4").Value = 0.038248577916091514
$ws5.Range("L4").Value = 0.00007953649306931412
$ws5.Range("N4").Value = 0.007761399678967117
$ws5.Range("B5").Value = 0.016967257085302008
$ws5.Range("C5").Value = 0.005201954480383295
$ws5.Range("D5").Value = 0.0003170465189210434
$ws5.Range("E5").Value = 0.02795693855715173
$ws5.Range("F5").Value = 0.0009352870840520827
$ws5.Range("G5").Value = 0.006254225654386913
$ws5.Range("I5").Value = 0.015915136903460986
$ws5.Range("J5").Value = 0.029060917574679154
$ws5.Range("L5").Value = 0.00008766704951244843
$ws5.Range("N5").Value = 0.003931401988305117
$ws5.Range("B6").Value = 0.016703966897314366
$ws5.Range("C6").Value = 0.006222025439678842
$ws5.Range("D6").Value = 0.00035719583796350293
$ws5.Range("E6").Value = 0.02745537400223108
$ws5.Range("F6").Value = 0.0008937618304399849
$ws5.Range("G6").Value = 0.009614610501094906
$ws5.Range("I6").Value = 0.014634778950769102
$ws5.Range("J6").Value = 0.025826827578832583
$ws5.Range("L6").Value = 0.00006601901645155402
$ws5.Range("N6").Value = 0.0058193001581905105
$ws5.Range("B7").Value = 0.0052171137279507975
$ws5.Range("C7").Value = 0.006208600626257491
$ws5.Range("D7").Value = 0.00005788260017667619
$ws5.Range("E7").Value = 0.019771611650357908
$ws5.Range("F7").Value = 0.0008362323706893776
$ws5.Range("G7").Value = 0.009326775729711183
$ws5.Range("I7").Value = 0.0010052958459337
$ws5.Range("J7").Value = 0.018482430287661514
$ws5.Range("L7").Value = 0.000007521653875384236
$ws5.Range("N7").Value = 0.002101287852381847
$ws5.Range("B8").Value = 0.0015396181338185248
$ws5.Range("C8").Value = 0.0020521867080031237
$ws5.Range("D8").Value = 0.000012272873485124463
$ws5.Range("E8").Value = 0.009190798436029969
$ws5.Range("F8").Value = 0.0008648012796832334
$ws5.Range("G8").Value = 0.004025217570975997
$ws5.Range("I8").Value = 0.00026626694117635953
$ws5.Range("J8").Value = 0.017317782005123915
$ws5.Range("L8").Value = 0.000003087382289844282
$ws5.Range("N8").Value = 0.0002269197164400777
$ws5.Range("B9").Value = 0.0006533066242901698
$ws5.Range("C9").Value = 0.0008933547199427522
$ws5.Range("E9").Value = 0.006196215011245403
$ws5.Range("F9").Value = 0.0008342055587528082
$ws5.Range("G9").Value = 0.003540187512134227
$ws5.Range("I9").Value = 0.00005969180520870553
$ws5.Range("J9").Value = 0.014329578028222857
$ws5.Range("L9").Value = 0.00000023025386894147772
$ws5.Range("N9").Value = 0.000013992253335041258
$ws5.Range("B10").Value = 0.00023190606786513415
$ws5.Range("C10").Value = 0.00031854513254508596
$ws5.Range("E10").Value = 0.0022477631024692
$ws5.Range("F10").Value = 0.0008523271378405172
$ws5.Range("G10").Value = 0.004046540458400295
$ws5.Range("J10").Value = 0.01408213048661392
